# Commit: "added man's search for meaning"
#
# Adds a new row to the "2024" reading-log table (Table3 on the "2024"
# worksheet) for the book "Man's Search for Meaning" by Victor Frankl.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# The book log is an Excel Table (ListObject) currently spanning A1:F5.
# Add a new row to it the way a user would -- this grows the table to
# A1:F6 and keeps the table/autofilter definitions in sync.
$lo = $ws.ListObjects.Item(1)
[void]$lo.ListRows.Add()

# Pick up the same cell formatting (borders / number formats) used by the
# row above it, just like Excel does when a table auto-extends.
$ws.Range("A5:F5").Copy()
$ws.Range("A6:F6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new book's details.
$ws.Range("A6").Value = "Man's Search for Meaning"
$ws.Range("B6").Value = "Victor Frankl"
$ws.Range("C6").Value = 45474   # Date Started:  7/1/2024
$ws.Range("D6").Value = 45519   # Date Finished: 8/15/2024
$ws.Range("E6").Value = "***"
$ws.Range("F6").Value = "good book overall"

# Leave the view the way it ended up after adding the row.
$excel.ActiveWindow.Zoom = 120
[void]$ws.Range("F8").Select()
